$wb = $excel.ActiveWorkbook

# Overview sheet: update the "Latest HO Xliff Generate Date" for the
# f68abb70 file row (row 3) to reflect the newly generated handback report.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-28 18:47:03"

# zh-cn sheet: update Correspond Handoff Datetime / Correspond Handback
# DateTime for the f68abb70 file row (row 3).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-28 18:46:57"
$wsZhCn.Range("K3").Value = "2016-08-28 18:47:25"

# de-de sheet: update Correspond Handoff Datetime / Correspond Handback
# DateTime for the f68abb70 file row (row 3).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-28 18:47:03"
$wsDeDe.Range("K3").Value = "2016-08-28 18:47:32"
